# Apply scheduled-runner market price & profit updates to each sheet's quest table.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()  # was -5016
$ws.Range("N51").ClearContents()  # was -8781.666499999999
$ws.Range("H108").Value = 25000
$ws.Range("J108").Value = 25000
$ws.Range("L108").Value = 25000
$ws.Range("N108").Value = -32680
$ws.Range("H116").Value = 6610.143
$ws.Range("I116").Value = 7593.5
$ws.Range("J116").Value = 5299
$ws.Range("K116").Value = 7593.5
$ws.Range("L116").Value = 5299
$ws.Range("M116").Value = -4151.5
$ws.Range("N116").Value = -12183
$ws.Range("H135").Value = 6458.8335
$ws.Range("I135").Value = 4990.5
$ws.Range("J135").Value = 9395.5
$ws.Range("K135").Value = 44914.5
$ws.Range("L135").Value = 84559.5
$ws.Range("M135").Value = -42379.5
$ws.Range("N135").Value = -89629.5
$ws.Range("H138").Value = 2590.7778
$ws.Range("I138").Value = 1849.091
$ws.Range("J138").Value = 3756.2856
$ws.Range("K138").Value = 5547.272999999999
$ws.Range("L138").Value = 11268.8568
$ws.Range("M138").Value = -407.2729999999992
$ws.Range("N138").Value = -21548.8568
$ws.Range("H141").Value = 4166.3335
$ws.Range("I141").Value = 5000
$ws.Range("K141").Value = 15000
$ws.Range("M141").Value = -9820

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4281543
$ws.Range("I32").Value = 4121633.8
$ws.Range("K32").Value = 4121633.8
$ws.Range("M32").Value = -4121346.8
$ws.Range("H63").Value = 13663.143
$ws.Range("I63").Value = 15773.667
$ws.Range("K63").Value = 15773.667
$ws.Range("M63").Value = -15087.667
$ws.Range("H66").Value = 13663.143
$ws.Range("I66").Value = 15773.667
$ws.Range("K66").Value = 78868.33499999999
$ws.Range("M66").Value = -75436.33499999999
$ws.Range("H74").Value = 1859.8889
$ws.Range("I74").Value = 1715.5
$ws.Range("J74").Value = 1975.4
$ws.Range("K74").Value = 1715.5
$ws.Range("L74").Value = 1975.4
$ws.Range("M74").Value = -841.5
$ws.Range("N74").Value = -3723.4
$ws.Range("H77").Value = 1859.8889
$ws.Range("I77").Value = 1715.5
$ws.Range("J77").Value = 1975.4
$ws.Range("K77").Value = 8577.5
$ws.Range("L77").Value = 9877
$ws.Range("M77").Value = -4209.5
$ws.Range("N77").Value = -18613
$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492
$ws.Range("H102").Value = 2593.375
$ws.Range("I102").Value = 2249.5715
$ws.Range("K102").Value = 2249.5715
$ws.Range("M102").Value = -627.5715
$ws.Range("H112").Value = 32000
$ws.Range("J112").Value = 32000
$ws.Range("L112").Value = 32000
$ws.Range("N112").Value = -34954
$ws.Range("H132").Value = 1186.75
$ws.Range("I132").Value = 1249.3334
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 3748.0002
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -1218.0002
$ws.Range("N132").Value = -8057

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()  # was -2060
$ws.Range("H94").Value = 724.5
$ws.Range("I94").Value = 950
$ws.Range("K94").Value = 950
$ws.Range("M94").Value = -499

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4072
$ws.Range("I16").Value = 4931
$ws.Range("J16").Value = 3749.875
$ws.Range("K16").Value = 4931
$ws.Range("L16").Value = 3749.875
$ws.Range("M16").Value = -4644
$ws.Range("N16").Value = -4323.875
$ws.Range("H92").Value = 29466
$ws.Range("J92").Value = 29466
$ws.Range("L92").Value = 29466
$ws.Range("N92").Value = -34458
$ws.Range("H94").Value = 189780.17
$ws.Range("J94").Value = 2224.5
$ws.Range("L94").Value = 2224.5
$ws.Range("N94").Value = -3126.5
$ws.Range("H113").Value = 4072
$ws.Range("I113").Value = 4931
$ws.Range("J113").Value = 3749.875
$ws.Range("K113").Value = 4931
$ws.Range("L113").Value = 3749.875
$ws.Range("M113").Value = -2761
$ws.Range("N113").Value = -8089.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 2000
$ws.Range("K110").Value = 6000
$ws.Range("M110").Value = -1910
$ws.Range("H132").Value = 2706.9092
$ws.Range("I132").Value = 2148
$ws.Range("J132").Value = 3377.6
$ws.Range("K132").Value = 19332
$ws.Range("L132").Value = 30398.4
$ws.Range("M132").Value = -16802
$ws.Range("N132").Value = -35458.39999999999
$ws.Range("H134").Value = 1899.3334
$ws.Range("I134").Value = 1899.3334
$ws.Range("K134").Value = 5698.0002
$ws.Range("M134").Value = -628.0002000000004
$ws.Range("H140").Value = 1036.6666
$ws.Range("I140").Value = 1036.6666
$ws.Range("K140").Value = 3109.9998
$ws.Range("M140").Value = 2070.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4003.5
$ws.Range("I70").Value = 4003.5
$ws.Range("K70").Value = 4003.5
$ws.Range("M70").Value = -3733.5
$ws.Range("H73").Value = 4003.5
$ws.Range("I73").Value = 4003.5
$ws.Range("K73").Value = 4003.5
$ws.Range("M73").Value = -3067.5
$ws.Range("H95").Value = 27997
$ws.Range("J95").Value = 27997
$ws.Range("L95").Value = 27997
$ws.Range("N95").Value = -33489

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8489.817999999999
$ws.Range("I7").Value = 7528.6
$ws.Range("K7").Value = 7528.6
$ws.Range("M7").Value = -7416.6
$ws.Range("H16").Value = 1682.1111
$ws.Range("I16").Value = 1672.6
$ws.Range("K16").Value = 1672.6
$ws.Range("M16").Value = -1502.6
$ws.Range("H22").Value = 2950
$ws.Range("I22").Value = 2950
$ws.Range("K22").Value = 2950
$ws.Range("M22").Value = -2655
$ws.Range("H27").Value = 2950
$ws.Range("I27").Value = 2950
$ws.Range("K27").Value = 2950
$ws.Range("M27").Value = -2843
$ws.Range("H40").Value = 5366.1665
$ws.Range("I40").Value = 5339.4
$ws.Range("K40").Value = 5339.4
$ws.Range("M40").Value = -5203.4
$ws.Range("H68").Value = 3199.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3199.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3199.5
$ws.Range("M68").ClearContents()  # was -1200.5
$ws.Range("N68").Value = -4697.5
$ws.Range("H71").Value = 3199.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3199.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 15997.5
$ws.Range("M71").ClearContents()  # was -6003.5
$ws.Range("N71").Value = -23485.5
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()  # was 123
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()  # was 2616
$ws.Range("H100").Value = 2270
$ws.Range("I100").Value = 2018.6
$ws.Range("J100").Value = 2898.5
$ws.Range("K100").Value = 2018.6
$ws.Range("L100").Value = 2898.5
$ws.Range("M100").Value = -1477.6
$ws.Range("N100").Value = -3980.5
$ws.Range("H122").Value = 7452.25
$ws.Range("I122").Value = 4816
$ws.Range("J122").Value = 8060.615
$ws.Range("K122").Value = 14448
$ws.Range("L122").Value = 24181.845
$ws.Range("M122").Value = -11998
$ws.Range("N122").Value = -29081.845
$ws.Range("H126").Value = 8489.817999999999
$ws.Range("I126").Value = 7528.6
$ws.Range("K126").Value = 22585.8
$ws.Range("M126").Value = -20115.8
$ws.Range("H132").Value = 3558.625
$ws.Range("I132").Value = 3367.5
$ws.Range("K132").Value = 10102.5
$ws.Range("M132").Value = -7572.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2501017.2
$ws.Range("I81").Value = 1355.6666
$ws.Range("K81").Value = 2711.3332
$ws.Range("M81").Value = -1650.3332
$ws.Range("H82").Value = 40000
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40766
$ws.Range("H84").Value = 2501017.2
$ws.Range("I84").Value = 1355.6666
$ws.Range("K84").Value = 13556.666
$ws.Range("M84").Value = -8252.666000000001
$ws.Range("H85").Value = 40000
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42652
$ws.Range("H126").Value = 2740.2727
$ws.Range("I126").Value = 2077.7144
$ws.Range("K126").Value = 6233.1432
$ws.Range("M126").Value = -3763.1432
